$d = $word.ActiveDocument
$r = $d.Content
$found = $r.Find.Execute("pour afficher les outils disponibles.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)

$cr = [char]13
$insertText = $cr + "Le 19 novembre 2014" + $cr + "J’ai débuté la scrollbar, mais je ne suis pas trop sûr de comment m’y prendre encore."
$r.InsertAfter($insertText)
